$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")

# Clear the region we are rewriting (old dimension was A2:D14)
$ws.Range("A2:E26").ClearContents()

# --- Header rows (moved from C2:D3 to D2:E3) ---
$ws.Range("D2").Value = 'Tue'
$ws.Range("E2").Value = 'Wed'
$ws.Range("D3").Value = 'Hours'
$ws.Range("E3").Value = 'Hours'

# --- New upper backlog block (rows 5-14) ---
$ws.Range("A5").Value = 'Validate Worksheet and Calc for SaskCrown Oil'
$ws.Range("B5").Value = 'L'
$ws.Range("D5").Value = 2
$ws.Range("A6").Value = 'Setup new project and new dirctory structure'
$ws.Range("B6").Value = 'K'
$ws.Range("A7").Value = 'Create a model module with tests and convensions (Where do our test go, nameing convension)'
$ws.Range("B7").Value = 'K'
$ws.Range("A8").Value = 'Implement Missing IOGR1995 Supplementary Royalty with tests'
$ws.Range("B8").Value = 'A'
$ws.Range("A9").Value = 'Write Tests and restructure calcSaskOilProvCrown'
$ws.Range("B9").Value = 'A'
$ws.Range("A10").Value = 'Write DB Load utilities for testing'
$ws.Range("A12").Value = 'Adrienne - Test Drivin Development, Python code'
$ws.Range("A13").Value = 'Konstantin - sqlite, calcs'
$ws.Range("A14").Value = 'Larry - Utilities help'

# --- Lower status block (rows 16-26) ---
$ws.Range("A16").Value = 'Get the code working as is: GUI, Calculation, Worksheet'
$ws.Range("B16").Value = 'K'
$ws.Range("C16").Value = 'Done'
$ws.Range("D16").Value = 1
$ws.Range("A17").Value = 'Add worksheet to GUI'
$ws.Range("B17").Value = 'K'
$ws.Range("C17").Value = 'Cancel'
$ws.Range("D17").Value = 3
$ws.Range("A18").Value = 'Demo and obtain feedback'
$ws.Range("B18").Value = 'L'
$ws.Range("C18").Value = 'Done'
$ws.Range("D18").Value = 14
$ws.Range("A20").Value = 'Adrienne Learn python test'
$ws.Range("B20").Value = 'A'
$ws.Range("C20").Value = 'Done'
$ws.Range("D20").Value = 7
$ws.Range("A21").Value = 'Get System running on Adrienne''s machine'
$ws.Range("B21").Value = 'K'
$ws.Range("C21").Value = 'Done'
$ws.Range("D21").Value = 3
$ws.Range("A22").Value = 'Get TDD running in Adrienne''s machine'
$ws.Range("B22").Value = 'K&A'
$ws.Range("A23").Value = 'Review and Learn Existing Calcs '
$ws.Range("B23").Value = 'A'
$ws.Range("C23").Value = 'Done'
$ws.Range("A25").Value = 'Create TDD Strategy  (Coverage, No Code without a test)'
$ws.Range("B25").Value = 'K'
$ws.Range("D25").Value = 10
$ws.Range("A26").Value = 'Create Data Access Strategy'
$ws.Range("B26").Value = 'L'
$ws.Range("D26").Value = 10

# --- Column widths: widen to cover new column C (Status) at the same width as B ---
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# --- Wrap text on column A, applied only to the populated cells (keeps gap rows empty) ---
$ws.Range("A5:A10").WrapText = $true
$ws.Range("A12:A14").WrapText = $true
$ws.Range("A16:A18").WrapText = $true
$ws.Range("A20:A23").WrapText = $true
$ws.Range("A25:A26").WrapText = $true

# --- Taller rows for the two-line wrapped descriptions ---
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 30

# --- Selection + page setup ---
$ws.Range("B12").Select()
$ws.PageSetup.Orientation = 1
